$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("BO2").Value = 0.7066429257392883
$ws.Range("BO3").Value = 0.8333324790000916
$ws.Range("BO4").Value = 0.6847591996192932
$ws.Range("BO5").Value = 0.7155214548110962
$ws.Range("BO6").Value = 0.7534450888633728
$ws.Range("BO7").Value = 0.6875506043434143
$ws.Range("BO8").Value = 0.7220718860626221
$ws.Range("BO9").Value = 0.7777664661407471
$ws.Range("BO10").Value = 0.8194712996482849
$ws.Range("BO11").Value = 0.738487184047699
$ws.Range("BO12").Value = 0.6883087754249573
$ws.Range("BO13").Value = 0.759907066822052
$ws.Range("BO14").Value = 0.7971763610839844
$ws.Range("BO15").Value = 0.6778339147567749
$ws.Range("BO16").Value = 0.8134178519248962
$ws.Range("BO17").Value = 0.6895470023155212
$ws.Range("BO18").Value = 0.7083225846290588
$ws.Range("BO19").Value = 0.03088539279997349
$ws.Range("BO20").Value = 0.6962499022483826
$ws.Range("BO21").Value = 0.2887141704559326
$ws.Range("BO22").Value = 0.4580791890621185
$ws.Range("BO23").Value = 0.7294402122497559
$ws.Range("BO24").Value = 0.6696252226829529
$ws.Range("BO25").Value = 0.6971911191940308
$ws.Range("BO26").Value = 0.7828356623649597
$ws.Range("BO27").Value = 0.6579923033714294
$ws.Range("BO28").Value = 0.6283414959907532
$ws.Range("BO29").Value = 0.6178016066551208
$ws.Range("BO30").Value = 0.6607226729393005
$ws.Range("BO31").Value = 0.773419201374054
$ws.Range("BO32").Value = 0.6031879186630249
$ws.Range("BO33").Value = 0.5983608365058899
$ws.Range("BO34").Value = 0.701208233833313
$ws.Range("BO35").Value = 0.714297354221344
$ws.Range("BO36").Value = 0.6608238816261292
$ws.Range("BO37").Value = 0.8220487236976624
$ws.Range("BO38").Value = 0.7736412882804871
$ws.Range("BO39").Value = 0.4430015981197357
$ws.Range("BO40").Value = 0.6461691856384277
$ws.Range("BO41").Value = 0.4895467758178711
$ws.Range("BO42").Value = 0.4848021864891052
$ws.Range("BO43").Value = 0.6732003092765808
$ws.Range("BO44").Value = 0.6191999316215515
$ws.Range("BO45").Value = 0.5700157284736633
$ws.Range("BO46").Value = 0.6737620830535889
$ws.Range("BO47").Value = 0.7631568312644958
$ws.Range("BO48").Value = 0.69390469789505
$ws.Range("BO49").Value = 0.7028796076774597
$ws.Range("BO50").Value = 0.6760109066963196
$ws.Range("BO51").Value = 0.6505436301231384
$ws.Range("BO52").Value = 0.1117662787437439
$ws.Range("BO53").Value = 0.1265599876642227
$ws.Range("BO54").Value = 0.00934301596134901
$ws.Range("BO55").Value = 0.6037424206733704
$ws.Range("BO56").Value = 0.7923555374145508
$ws.Range("BO57").Value = 0.6299933195114136
$ws.Range("BO58").Value = 0.8128296732902527
$ws.Range("BO59").Value = 0.7638512253761292
$ws.Range("BO60").Value = 0.8637326955795288
$ws.Range("BO61").Value = 0.7564572691917419
$ws.Range("BO62").Value = 0.7503859996795654
$ws.Range("BO63").Value = 0.6832898259162903
$ws.Range("BO64").Value = 0.6595971584320068
$ws.Range("BO65").Value = 0.6625831127166748
$ws.Range("BO66").Value = 0.08363591134548187
$ws.Range("BO67").Value = 0.6920419931411743
$ws.Range("BO68").Value = 0.04806749150156975
$ws.Range("BO69").Value = 0.6870089173316956
$ws.Range("BO70").Value = 0.6720928549766541
$ws.Range("BO71").Value = 0.7154821157455444
$ws.Range("BO72").Value = 0.719618558883667
$ws.Range("BO73").Value = 0.8204642534255981
$ws.Range("BO74").Value = 0.8189181089401245
$ws.Range("BO75").Value = 0.7602909803390503
$ws.Range("BO76").Value = 0.7170021533966064
